$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Saudi Arabia" row from each of the three groups (Global, High-income,
# International). Row 13 is Saudi Arabia in the first (Global) group; after it is
# removed the following groups shift up by one row, so the next Saudi Arabia row
# (originally row 26) is now row 25, and the last one (originally row 39) is now row 37.
$ws.Rows(13).Delete()
$ws.Rows(25).Delete()
$ws.Rows(37).Delete()

# Re-run (updated) values for the "<b>All</b>" aggregate row and the "Russia" row in
# each of the three groups.

# Global group
$ws.Range("B2").Value = 73.6710361597934
$ws.Range("C2").Value = 72.3033300855294
$ws.Range("D2").Value = 75.0387422340573

$ws.Range("B12").Value = 77.5532714938846
$ws.Range("C12").Value = 73.1395709765103
$ws.Range("D12").Value = 81.966972011259

# High-income group
$ws.Range("B14").Value = 68.7235111211832
$ws.Range("C14").Value = 67.2867144309077
$ws.Range("D14").Value = 70.1603078114586

$ws.Range("B24").Value = 69.8291634341244
$ws.Range("C24").Value = 64.9837035538317
$ws.Range("D24").Value = 74.6746233144171

# International group
$ws.Range("B26").Value = 67.4325771642063
$ws.Range("C26").Value = 65.9835497503382
$ws.Range("D26").Value = 68.8816045780744

$ws.Range("B36").Value = 73.9042764455841
$ws.Range("C36").Value = 69.0236003408544
$ws.Range("D36").Value = 78.7849525503137
